# Atualizar dados de faturamento diario das lojas Bibi (atualizei dados da bibi e add)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Linha 2 - Bibi Cell Mundi
$ws.Range("S2").Value = 44502.09
$ws.Range("AG2").Value = 233545.67

# Linha 3 - Bibi Cell Vieiralves
$ws.Range("Q3").Value = 2253
$ws.Range("S3").Value = 3548
$ws.Range("AG3").Value = 84387.3

# Linha 4 - Bibi Cell Ponta Negra
$ws.Range("S4").Value = 2255.01
$ws.Range("AG4").Value = 55367.33

# Linha 5 - Bibi Cell Manauara
$ws.Range("S5").Value = 1209
$ws.Range("AG5").Value = 49004.05

# Linha 6 - total
$ws.Range("Q6").Value = 21142.5
$ws.Range("S6").Value = 51514.1
$ws.Range("AG6").Value = 422304.35

$wb.Save()
